# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the two sets of header columns from the generic "_old"/"_new"
# suffixes to the concrete format-version suffixes "_FV2410"/"_FV2504",
# wraps the used range in an Excel Table ("Table1"), and freezes the
# header row (row 1) with a pane split.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells (row 1) ------------------------------------
# Columns A:K were suffixed "_old" -> now "_FV2410"
# Columns L:U were suffixed "_new" -> now "_FV2504"
$headerRange = $ws.Range("A1:U1")
for ($c = 1; $c -le $headerRange.Columns.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value()
    if ($null -ne $text -and $text -ne "") {
        $text = $text -replace "_old$", "_FV2410"
        $text = $text -replace "_new$", "_FV2504"
        $cell.Value = $text
    }
}

# --- 2) Wrap the data range in an Excel Table ---------------------------
$usedRange = $ws.Range("A1:U55")
$table = $ws.ListObjects.Add(0, $usedRange, $null, 1)
$table.Name = "Table1"
# Drop the auto-assigned default style name so the exported table keeps a
# bare <tableStyleInfo> (matching the source file, which wasn't styled).
$table.TableStyle = ""

# --- 3) Freeze the header row (pane split) ------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
